$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 222, shifting existing rows 222-264 down to 223-265
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row 222 with the new data record
$ws.Cells.Item(222, 1).Value = 9
$ws.Cells.Item(222, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(222, 3).Value = "Metropolitana"
$ws.Cells.Item(222, 4).Value = 44694
$ws.Cells.Item(222, 5).Value = 13
$ws.Cells.Item(222, 6).Value = 300000001
$ws.Cells.Item(222, 7).Value = "Rabanito"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 6100
$ws.Cells.Item(222, 11).Value = 2500
$ws.Cells.Item(222, 12).Value = 3000
$ws.Cells.Item(222, 13).Value = 2750
$ws.Cells.Item(222, 14).Value = "`$/cien unidades (volumen en unidades)"
$ws.Cells.Item(222, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(222, 16).Value = 28
$ws.Cells.Item(222, 17).Value = 100
$ws.Cells.Item(222, 18).Value = "Hortaliza"

# Apply the same date style as other date cells in column D (style index 2 -> numFmtId 165)
$ws.Cells.Item(222, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
